$wb = $excel.ActiveWorkbook

$sheetNames = @("Stations_Mean", "Stations_Std", "Stations_CV")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B2").Value = "chl (ug/l)"
    $ws.Range("C2").Value = "turbidity (NTU)"
}
